$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Restructure the "Sheet1" test-case data: insert a new "Reports" block
#    (1 blank separator row + 9 data rows) right after the "Login" block,
#    and remove the obsolete "Product" block (1 blank separator + 2 data rows).
# ---------------------------------------------------------------------------

# Insert 10 new rows right after row 5 (the end of the "Login" block). This
# pushes "Category" (and everything after it) down from row 7 to row 17,
# and copies the D-column "Pass" formatting down into the new rows.
$ws.Range("A6:A15").EntireRow.Insert()

# Give the new "Reports" label cell (A7) the same bold category-label style
# used by the other section headers (copy format only from "Category", A17).
$ws.Range("A17").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Fill in the new "Reports" section (rows 7-15).
$ws.Range("A7").Value() = "Reports"

$ws.Range("B7").Value() = " View Review and Rating Report"
$ws.Range("C7").Value() = "Review and Rating report is displayed."
$ws.Range("D7").Value() = "Pass"

$ws.Range("B8").Value() = " View  Products Report"
$ws.Range("C8").Value() = "Products report is displayed."
$ws.Range("D8").Value() = "Pass"

$ws.Range("B9").Value() = " View  Traders Report"
$ws.Range("C9").Value() = "Traders report is displayed."
$ws.Range("D9").Value() = "Pass"

$ws.Range("B10").Value() = " View  Daily Orders Report"
$ws.Range("C10").Value() = "Daily orders report is displayed."
$ws.Range("D10").Value() = "Pass"

$ws.Range("B11").Value() = " View Weekly Order Report"
$ws.Range("C11").Value() = "Weekly orders Report id displayed."
$ws.Range("D11").Value() = "Pass"

$ws.Range("B12").Value() = " View  Monthly Order Report"
$ws.Range("C12").Value() = "Monthly Orders Report is displayed."
$ws.Range("D12").Value() = "Pass"

$ws.Range("B13").Value() = " View Daily Payment Report"
$ws.Range("C13").Value() = "Daily Payment Report is displayed."
$ws.Range("D13").Value() = "Pass"

$ws.Range("B14").Value() = " View Weekly Payment Report"
$ws.Range("C14").Value() = "Weekly Payment Report is displayed."
$ws.Range("D14").Value() = "Pass"

$ws.Range("B15").Value() = " View Monthly Payment Report"
$ws.Range("C15").Value() = "Monthly Payment Report is displayed."
$ws.Range("D15").Value() = "Pass"

# Remove the obsolete "Product" block. After the insert above it now lives
# at rows 31 (blank separator), 32 ("Activate Product") and 33
# ("Deactivate Product") - delete all three so "Profile" moves straight up.
$ws.Range("A31:A33").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Add a new blank worksheet ("Sheet4") in front of the other sheets and
#    make "Sheet1" (now holding the restructured data) the active tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item("Sheet1").Activate()
